# Apply targeted cell updates to Sheet1 to reflect the corrected
# classification values ("task 1 notebook update with few shot").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$updates = @{
    "C2"  = "Low"
    "A7"  = "None"
    "A8"  = "Gampaha"
    "A9"  = "None"
    "A10" = "None"
    "A11" = "None"
    "A13" = "None"
    "A14" = "Colombo"
    "A15" = "Colombo"
    "A16" = "None"
    "A19" = "None"
    "A20" = "None"
    "A21" = "None"
    "A22" = "None"
    "A23" = "None"
    "A24" = "None"
    "A30" = "None"
    "C32" = "Low"
    "A33" = "None"
    "B33" = "Other"
    "A34" = "None"
    "B34" = "Other"
    "A37" = "None"
    "A39" = "None"
    "A41" = "None"
    "B41" = "Other"
    "B42" = "Supply"
    "C42" = "Low"
    "A43" = "None"
    "B43" = "Rescue"
    "C43" = "High"
    "A44" = "None"
    "A45" = "Kandy"
    "A47" = "None"
    "B47" = "Other"
    "A48" = "None"
    "A49" = "None"
    "A50" = "None"
    "B50" = "Info"
    "A51" = "None"
    "B51" = "Other"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
